# data-system-share1_v2: update title slide credit line and drop the
# trailing "Any question?" closing slide (and its notes page).

$p = $ppt.ActivePresentation

# --- 1. Title slide: drop the "主讲人：运维 · " lead-in before the
#        presenter's name and center the remaining line. -----------------
$titleSlide = $p.Slides.Item(1)
$creditShape = $titleSlide.Shapes.Item(2)
$creditRange = $creditShape.TextFrame.TextRange

# Remove the first 9 characters ("主讲人：运维 " + "· ") while leaving the
# "朱卫中" runs (and their formatting) untouched.
$lead = $creditRange.Characters(1, 9)
$lead.Text = ""

# Center-align the now-shorter paragraph.
$creditRange.ParagraphFormat.Alignment = 2

# --- 2. Remove the final "Any question?" slide (and its notes page
#        cascades automatically with it). ---------------------------------
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()
